# Apply scheduled market-price / profit recalculations to the Sheets workbook.
# Each worksheet corresponds to a crafting class (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR);
# cells H:N hold price/profit figures that are refreshed by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 233.5
$ws.Range("J2").Value = 351
$ws.Range("L2").Value = 351
$ws.Range("N2").Value = -577
# Row 15
$ws.Range("H15").Value = 2368.6575
$ws.Range("I15").Value = 2368.6575
$ws.Range("K15").Value = 7105.9725
$ws.Range("M15").Value = -6936.9725
# Row 19
$ws.Range("H19").Value = 10987.35
$ws.Range("I19").Value = 874.7778
$ws.Range("J19").Value = 19261.273
$ws.Range("K19").Value = 874.7778
$ws.Range("L19").Value = 19261.273
$ws.Range("M19").Value = -699.7778
$ws.Range("N19").Value = -19611.273
# Row 33
$ws.Range("H33").Value = 1714.9375
$ws.Range("I33").Value = 1549
$ws.Range("J33").Value = 2080
$ws.Range("K33").Value = 1549
$ws.Range("L33").Value = 2080
$ws.Range("M33").Value = -1320
$ws.Range("N33").Value = -2538
# Row 40
$ws.Range("H40").Value = 1326.5555
$ws.Range("I40").Value = 1268
$ws.Range("J40").Value = 1399.75
$ws.Range("K40").Value = 1268
$ws.Range("L40").Value = 1399.75
$ws.Range("M40").Value = -1093
$ws.Range("N40").Value = -1749.75
# Row 43
$ws.Range("H43").Value = 1768.1428
$ws.Range("J43").Value = 2202.6667
$ws.Range("L43").Value = 2202.6667
$ws.Range("N43").Value = -2340.6667
# Row 64
$ws.Range("H64").Value = 3138.889
$ws.Range("I64").Value = 3147.0588
$ws.Range("K64").Value = 3147.0588
$ws.Range("M64").Value = -2899.0588
# Row 67
$ws.Range("H67").Value = 3138.889
$ws.Range("I67").Value = 3147.0588
$ws.Range("K67").Value = 3147.0588
$ws.Range("M67").Value = -2289.0588
# Row 127
$ws.Range("H127").Value = 1311.6428
$ws.Range("I127").Value = 1100
$ws.Range("J127").Value = 1369.3636
$ws.Range("K127").Value = 3300
$ws.Range("L127").Value = 4108.0908
$ws.Range("M127").Value = 1660
$ws.Range("N127").Value = -14028.0908
# Row 138
$ws.Range("H138").Value = 136328.9
$ws.Range("J138").Value = 183322.25
$ws.Range("L138").Value = 549966.75
$ws.Range("N138").Value = -560246.75

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3285.7144
$ws.Range("I61").Value = 2333.3333
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2333.3333
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2121.3333
$ws.Range("N61").Value = -4424
# Row 102
$ws.Range("H102").Value = 1500
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 1500
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = $null
$ws.Range("N102").Value = -4744
# Row 122
$ws.Range("H122").Value = 2553.8
$ws.Range("I122").Value = 1904
$ws.Range("J122").Value = 3528.5
$ws.Range("K122").Value = 5712
$ws.Range("L122").Value = 10585.5
$ws.Range("M122").Value = -3262
$ws.Range("N122").Value = -15485.5
# Row 123
$ws.Range("H123").Value = 34428
$ws.Range("J123").Value = 34428
$ws.Range("L123").Value = 34428
$ws.Range("N123").Value = -44228
# Row 132
$ws.Range("H132").Value = 6737.52
$ws.Range("I132").Value = 7330.7144
$ws.Range("J132").Value = 5982.5454
$ws.Range("K132").Value = 21992.1432
$ws.Range("L132").Value = 17947.6362
$ws.Range("M132").Value = -19462.1432
$ws.Range("N132").Value = -23007.6362
# Row 136
$ws.Range("H136").Value = 3285.7144
$ws.Range("I136").Value = 2333.3333
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 6999.999899999999
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -4449.999899999999
$ws.Range("N136").Value = -17100

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 2266.2964
$ws.Range("I22").Value = 2420.8
$ws.Range("K22").Value = 2420.8
$ws.Range("M22").Value = -2247.8

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 367.5
$ws.Range("I22").Value = 277.5
$ws.Range("K22").Value = 277.5
$ws.Range("M22").Value = 72.5
# Row 31
$ws.Range("H31").Value = 2530.4517
$ws.Range("I31").Value = 1052.963
$ws.Range("J31").Value = 12503.5
$ws.Range("K31").Value = 1052.963
$ws.Range("L31").Value = 12503.5
$ws.Range("M31").Value = -757.963
$ws.Range("N31").Value = -13093.5
# Row 34
$ws.Range("H34").Value = 2530.4517
$ws.Range("I34").Value = 1052.963
$ws.Range("J34").Value = 12503.5
$ws.Range("K34").Value = 1052.963
$ws.Range("L34").Value = 12503.5
$ws.Range("M34").Value = -850.963
$ws.Range("N34").Value = -12907.5
# Row 58
$ws.Range("H58").Value = 1519.4
$ws.Range("I58").Value = 1466.6666
$ws.Range("J58").Value = 1528.7059
$ws.Range("K58").Value = 1466.6666
$ws.Range("L58").Value = 1528.7059
$ws.Range("M58").Value = -1263.6666
$ws.Range("N58").Value = -1934.7059
# Row 62
$ws.Range("H62").Value = 3689
$ws.Range("I62").Value = 2961
$ws.Range("J62").Value = 4599
$ws.Range("K62").Value = 2961
$ws.Range("L62").Value = 4599
$ws.Range("N62").Value = -5847
$ws.Range("M62").Value = -2337
# Row 65
$ws.Range("H65").Value = 3689
$ws.Range("I65").Value = 2961
$ws.Range("J65").Value = 4599
$ws.Range("K65").Value = 14805
$ws.Range("L65").Value = 22995
$ws.Range("N65").Value = -29235
$ws.Range("M65").Value = -11685
# Row 136
$ws.Range("H136").Value = 1519.4
$ws.Range("I136").Value = 1466.6666
$ws.Range("J136").Value = 1528.7059
$ws.Range("K136").Value = 4399.9998
$ws.Range("L136").Value = 4586.1177
$ws.Range("M136").Value = -1849.9998
$ws.Range("N136").Value = -9686.117699999999

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1237.45
$ws.Range("I113").Value = 699.9
$ws.Range("J113").Value = 1775
$ws.Range("K113").Value = 2099.7
$ws.Range("L113").Value = 5325
$ws.Range("M113").Value = 70.30000000000018
$ws.Range("N113").Value = -9665
# Row 125
$ws.Range("H125").Value = 2108.6
$ws.Range("J125").Value = 2108.6
$ws.Range("L125").Value = 6325.799999999999
$ws.Range("N125").Value = -16165.8
# Row 131
$ws.Range("H131").Value = 852.25714
$ws.Range("I131").Value = 222.5
$ws.Range("J131").Value = 1038.8518
$ws.Range("K131").Value = 667.5
$ws.Range("L131").Value = 3116.5554
$ws.Range("M131").Value = 4372.5
$ws.Range("N131").Value = -13196.5554
# Row 138
$ws.Range("H138").Value = 2687.1562
$ws.Range("I138").Value = 2358.3333
$ws.Range("J138").Value = 2763.0386
$ws.Range("K138").Value = 7074.999899999999
$ws.Range("L138").Value = 8289.1158
$ws.Range("M138").Value = -1934.999899999999
$ws.Range("N138").Value = -18569.1158
# Row 140
$ws.Range("H140").Value = 1822.5714
$ws.Range("I140").Value = 1269.5714
$ws.Range("J140").Value = 2928.5715
$ws.Range("K140").Value = 3808.7142
$ws.Range("L140").Value = 8785.7145
$ws.Range("M140").Value = 1371.2858
$ws.Range("N140").Value = -19145.7145

$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 43801.5
$ws.Range("J22").Value = 43801.5
$ws.Range("L22").Value = 43801.5
$ws.Range("N22").Value = -44859.5
# Row 25
$ws.Range("H25").Value = 80009
$ws.Range("J25").Value = 80009
$ws.Range("L25").Value = 80009
$ws.Range("N25").Value = -81067
# Row 122
$ws.Range("H122").Value = 5273.76
$ws.Range("J122").Value = 5523.4736
$ws.Range("L122").Value = 16570.4208
$ws.Range("N122").Value = -21470.4208

$ws = $wb.Worksheets.Item("LTW")
# Row 51
$ws.Range("H51").Value = 14000
$ws.Range("J51").Value = 14000
$ws.Range("L51").Value = 14000
$ws.Range("N51").Value = -14956
# Row 74
$ws.Range("H74").Value = 26499.5
$ws.Range("I74").Value = 19999
$ws.Range("J74").Value = 33000
$ws.Range("K74").Value = 19999
$ws.Range("L74").Value = 33000
$ws.Range("N74").Value = -34996
$ws.Range("M74").Value = -19001
# Row 77
$ws.Range("H77").Value = 26499.5
$ws.Range("I77").Value = 19999
$ws.Range("J77").Value = 33000
$ws.Range("K77").Value = 59997
$ws.Range("L77").Value = 99000
$ws.Range("N77").Value = -108984
$ws.Range("M77").Value = -55005
# Row 122
$ws.Range("H122").Value = 5300.4
$ws.Range("I122").Value = 4857.7144
$ws.Range("K122").Value = 14573.1432
$ws.Range("M122").Value = -12123.1432
# Row 132
$ws.Range("H132").Value = 4596.8965
$ws.Range("I132").Value = 3994.6667
$ws.Range("J132").Value = 5242.143
$ws.Range("K132").Value = 11984.0001
$ws.Range("L132").Value = 15726.429
$ws.Range("M132").Value = -9454.000100000001
$ws.Range("N132").Value = -20786.429
# Row 136
$ws.Range("H136").Value = 5954335
$ws.Range("I136").Value = 2253.1333
$ws.Range("J136").Value = 12822122
$ws.Range("K136").Value = 6759.3999
$ws.Range("L136").Value = 38466366
$ws.Range("M136").Value = -4209.3999
$ws.Range("N136").Value = -38471466

$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872
# Row 78
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360
# Row 123
$ws.Range("H123").Value = 30024.908
$ws.Range("J123").Value = 30024.908
$ws.Range("L123").Value = 30024.908
$ws.Range("N123").Value = -39824.908
# Row 136
$ws.Range("H136").Value = 3071.4827
$ws.Range("I136").Value = 2485.8125
$ws.Range("K136").Value = 7457.4375
$ws.Range("M136").Value = -4907.4375
